$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("UserList")

# Employee Status column (G) now stores text "t"/"f" instead of numeric 1/0
$ws1.Range("G2").Value = "t"
$ws1.Range("G2").Style = "Normal"
$ws1.Range("G3").Value = "t"
$ws1.Range("G4").Value = "t"
$ws1.Range("G5").Value = "t"
$ws1.Range("G6").Value = "f"
$ws1.Range("G7").Value = "f"
$ws1.Range("G8").Value = "f"
$ws1.Range("G9").Value = "f"

# Card # column (E) gets an explicit integer number format; new empty
# formatted cells appear down through row 9
$ws1.Range("E2:E9").NumberFormat = "0"

# Column E widens to fit the new formatting / content
$ws1.Range("E1").ColumnWidth = 11.74

# Make UserList the active/selected sheet with F6 selected
# (RequestList keeps its own F3 selection, it's just no longer the active tab)
$ws1.Activate() | Out-Null
$ws1.Range("F6").Select() | Out-Null
